$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the file-name related cells (C2, D2, E2, F2, E3, F3) to reflect the
# new test data file names used by the script.
$ws.Range("C2").Value = "extraction-template.xlsx"
$ws.Range("D2").Value = "\ExtractionTemplate\extraction-template.xlsx"
$ws.Range("E2").Value = "IC AML - Success Case Sheet.xlsx"
$ws.Range("F2").Value = "\ExtractionTemplate\ImportPublications\IC AML - Success Case Sheet.xlsx"
$ws.Range("E3").Value = "ICER - Failure Case Sheet.xlsx"
$ws.Range("F3").Value = "\ExtractionTemplate\ImportPublications\ICER - Failure Case Sheet.xlsx"

# Column D previously had its width auto-fit to its (now shortened) content;
# set the resulting best-fit width directly.
$ws.Columns.Item(4).ColumnWidth = 39.21875

# Move the active cell selection to C2, matching the saved view state.
$ws.Range("C2").Select()
